# fix: inversión inicial $3.500.000 y dividendo $0 en meses de gracia
#  - Inversión inicial actualizada de $3.200.000 a $3.500.000 CLP (Parámetros!B18)
#  - Dividendo / Arriendo (Flujo de Caja Mensual, fila 14) se aplica con IF:
#    $0 durante los meses de gracia ('Parámetros'!$B$5), y el valor de
#    'Parámetros'!$B$15 desde el mes siguiente en adelante.

$wb = $excel.ActiveWorkbook

# 1) Inversión inicial: Parámetros!B18 3.200.000 -> 3.500.000
$paramWs = $wb.Worksheets.Item("Parámetros")
$paramWs.Range("B18").Value = 3500000

# 2) Dividendo / Arriendo con meses de gracia: Flujo de Caja Mensual!B14:AK14
$flujoWs = $wb.Worksheets.Item("Flujo de Caja Mensual")

for ($col = 2; $col -le 37; $col++) {
    $mes = $col - 1
    $celda = $flujoWs.Cells.Item(14, $col)
    $celda.Formula = "=IF(" + $mes + "<='Parámetros'!`$B`$5,0,'Parámetros'!`$B`$15)"
}
